# #5: insurance, claim, debt, investment done
#
# Extends the "保險" (Insurance) and "債務" (Debt) sheets with the
# standard trailing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index [+ species/debtor/
# owner/total/register_date/register_reason for the debt sheet]) that all
# the other property sheets in this workbook already carry, and fixes the
# header row (row 1) on both sheets so it holds real field-name labels
# instead of a stray copy of row 2's data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "保險" (Insurance)
# ---------------------------------------------------------------------
$ins = $wb.Worksheets.Item("保險")

# Header row (row 1): proper field names for every column B..K
$ins.Cells.Item(1, 2).Value  = "company"
$ins.Cells.Item(1, 3).Value  = "name"
$ins.Cells.Item(1, 4).Value  = "owner"
$ins.Cells.Item(1, 5).Value  = "property_category"
$ins.Cells.Item(1, 6).Value  = "category"
$ins.Cells.Item(1, 7).Value  = "date"
$ins.Cells.Item(1, 8).Value  = "legislator_name"
$ins.Cells.Item(1, 9).Value  = "legislator_id"
$ins.Cells.Item(1, 10).Value = "source_file"
$ins.Cells.Item(1, 11).Value = "index"

# Data rows 2..12: columns B (company) / C (name) / D (owner) already
# hold correct values -- append the shared trailing metadata columns
# E..K, and clear the old stray E-column remark (register_reason) that
# doesn't belong to this schema any more.
$insRows = @(
    @{ Row = 2;  A = 104; B = "國泰人壽";     C = "新鍾情終身壽險";        D = "劉櫂豪" },
    @{ Row = 3;  A = 105; B = "國泰人壽";     C = "鍾意終身";              D = "劉櫂豪" },
    @{ Row = 4;  A = 106; B = "國泰人壽";     C = "創世紀變額萬能壽險(丁型）"; D = "劉櫂豪" },
    @{ Row = 5;  A = 107; B = "國泰人壽";     C = "安康住院醫療終身麵";    D = "林子煊" },
    @{ Row = 6;  A = 108; B = "國泰人壽";     C = "全福101終身";           D = "林子煊" },
    @{ Row = 7;  A = 109; B = "國泰人壽";     C = "得意還本終身";          D = "劉櫂豪" },
    @{ Row = 8;  A = 110; B = "國泰人壽";     C = "富貴年年終身";          D = "劉櫂豪" },
    @{ Row = 9;  A = 111; B = "富邦人壽";     C = "終身壽險";              D = "劉櫂豪" },
    @{ Row = 10; A = 112; B = "中國人壽";     C = "金享受終身壽險";        D = "劉櫂豪" },
    @{ Row = 11; A = 113; B = "國際紐約人壽"; C = "永安終身壽險";          D = "林子煊" },
    @{ Row = 12; A = 114; B = "富邦人壽";     C = "增美利外幣終身壽險";    D = "林子煊" }
)

foreach ($r in $insRows) {
    $row = $r.Row
    $ins.Cells.Item($row, 1).Value  = $r.A
    $ins.Cells.Item($row, 2).Value  = $r.B
    $ins.Cells.Item($row, 3).Value  = $r.C
    $ins.Cells.Item($row, 4).Value  = $r.D
    $ins.Cells.Item($row, 5).Value  = "insurance"
    $ins.Cells.Item($row, 6).Value  = "normal"
    $ins.Cells.Item($row, 7).Value  = "2012-04-30"
    $ins.Cells.Item($row, 8).Value  = "劉櫂豪"
    $ins.Cells.Item($row, 9).Value  = 1762
    $ins.Cells.Item($row, 10).Value = "tmpba991"
    $ins.Cells.Item($row, 11).Value = $r.A
}

# ---------------------------------------------------------------------
# Sheet "債務" (Debt)
# ---------------------------------------------------------------------
$debt = $wb.Worksheets.Item("債務")

# Header row (row 1): proper field names for every column B..N
$debt.Cells.Item(1, 2).Value  = "species"
$debt.Cells.Item(1, 3).Value  = "debtor"
$debt.Cells.Item(1, 4).Value  = "owner"
$debt.Cells.Item(1, 5).Value  = "total"
$debt.Cells.Item(1, 6).Value  = "register_date"
$debt.Cells.Item(1, 7).Value  = "register_reason"
$debt.Cells.Item(1, 8).Value  = "property_category"
$debt.Cells.Item(1, 9).Value  = "category"
$debt.Cells.Item(1, 10).Value = "date"
$debt.Cells.Item(1, 11).Value = "legislator_name"
$debt.Cells.Item(1, 12).Value = "legislator_id"
$debt.Cells.Item(1, 13).Value = "source_file"
$debt.Cells.Item(1, 14).Value = "index"

# Data row 2: keep/confirm existing species/debtor/owner/total/
# register_date/register_reason values and append the trailing
# metadata columns H..N.
$debt.Cells.Item(2, 1).Value  = 124
$debt.Cells.Item(2, 2).Value  = "房屋貸款"
$debt.Cells.Item(2, 3).Value  = "林子煊"
$debt.Cells.Item(2, 4).Value  = "新光商業銀行高雄七賢分行高雄市新興區七賢一路"
$debt.Cells.Item(2, 5).Value  = 2342159
$debt.Cells.Item(2, 6).Value  = "95年08月28日"
$debt.Cells.Item(2, 7).Value  = "購買房屋抵押貸款"
$debt.Cells.Item(2, 8).Value  = "debt"
$debt.Cells.Item(2, 9).Value  = "normal"
$debt.Cells.Item(2, 10).Value = "2012-04-30"
$debt.Cells.Item(2, 11).Value = "劉櫂豪"
$debt.Cells.Item(2, 12).Value = 1762
$debt.Cells.Item(2, 13).Value = "tmpba991"
$debt.Cells.Item(2, 14).Value = 124
